$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/updated figures for the 2020-08-03 data refresh.
# Each row's nombre_aides (C) and montant_total (D) are updated in place,
# keeping the cells stored as text (matching the source data format).

$c = $ws.Range("C5")
$c.NumberFormat = "@"
$c.Value = "121"
$c.ClearFormats()
$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "333886.40"
$d.ClearFormats()

$c = $ws.Range("C7")
$c.NumberFormat = "@"
$c.Value = "49"
$c.ClearFormats()
$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "121000.00"
$d.ClearFormats()

$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Value = "624"
$c.ClearFormats()
$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "1966440.56"
$d.ClearFormats()

$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Value = "31"
$c.ClearFormats()
$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "73270.00"
$d.ClearFormats()

$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Value = "127"
$c.ClearFormats()
$d = $ws.Range("D12")
$d.NumberFormat = "@"
$d.Value = "338000.00"
$d.ClearFormats()

$c = $ws.Range("C13")
$c.NumberFormat = "@"
$c.Value = "66"
$c.ClearFormats()
$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "159400.00"
$d.ClearFormats()

$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = "80"
$c.ClearFormats()
$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "200688.98"
$d.ClearFormats()

$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Value = "91"
$c.ClearFormats()
$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "334861.32"
$d.ClearFormats()

$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Value = "162"
$c.ClearFormats()
$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = "364100.00"
$d.ClearFormats()

$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = "28"
$c.ClearFormats()
$d = $ws.Range("D33")
$d.NumberFormat = "@"
$d.Value = "88000.00"
$d.ClearFormats()

$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "335"
$c.ClearFormats()
$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = "1271388.03"
$d.ClearFormats()

$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = "41"
$c.ClearFormats()
$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = "125000.00"
$d.ClearFormats()

$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = "74"
$c.ClearFormats()
$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = "167893.00"
$d.ClearFormats()

$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = "14"
$c.ClearFormats()
$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "48500.00"
$d.ClearFormats()

$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = "106"
$c.ClearFormats()
$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "273468.33"
$d.ClearFormats()

$c = $ws.Range("C79")
$c.NumberFormat = "@"
$c.Value = "394"
$c.ClearFormats()
$d = $ws.Range("D79")
$d.NumberFormat = "@"
$d.Value = "1458004.70"
$d.ClearFormats()

$c = $ws.Range("C88")
$c.NumberFormat = "@"
$c.Value = "94"
$c.ClearFormats()
$d = $ws.Range("D88")
$d.NumberFormat = "@"
$d.Value = "227510.00"
$d.ClearFormats()

$c = $ws.Range("C89")
$c.NumberFormat = "@"
$c.Value = "14"
$c.ClearFormats()
$d = $ws.Range("D89")
$d.NumberFormat = "@"
$d.Value = "31555.00"
$d.ClearFormats()

$c = $ws.Range("C90")
$c.NumberFormat = "@"
$c.Value = "53"
$c.ClearFormats()
$d = $ws.Range("D90")
$d.NumberFormat = "@"
$d.Value = "145495.14"
$d.ClearFormats()

$c = $ws.Range("C91")
$c.NumberFormat = "@"
$c.Value = "46"
$c.ClearFormats()
$d = $ws.Range("D91")
$d.NumberFormat = "@"
$d.Value = "125000.00"
$d.ClearFormats()

$c = $ws.Range("C92")
$c.NumberFormat = "@"
$c.Value = "101"
$c.ClearFormats()
$d = $ws.Range("D92")
$d.NumberFormat = "@"
$d.Value = "243405.00"
$d.ClearFormats()

$c = $ws.Range("C93")
$c.NumberFormat = "@"
$c.Value = "21"
$c.ClearFormats()
$d = $ws.Range("D93")
$d.NumberFormat = "@"
$d.Value = "43500.00"
$d.ClearFormats()

$c = $ws.Range("C94")
$c.NumberFormat = "@"
$c.Value = "123"
$c.ClearFormats()
$d = $ws.Range("D94")
$d.NumberFormat = "@"
$d.Value = "356902.00"
$d.ClearFormats()

$c = $ws.Range("C95")
$c.NumberFormat = "@"
$c.Value = "10"
$c.ClearFormats()
$d = $ws.Range("D95")
$d.NumberFormat = "@"
$d.Value = "24500.00"
$d.ClearFormats()

$c = $ws.Range("C96")
$c.NumberFormat = "@"
$c.Value = "7"
$c.ClearFormats()
$d = $ws.Range("D96")
$d.NumberFormat = "@"
$d.Value = "15500.00"
$d.ClearFormats()

$c = $ws.Range("C97")
$c.NumberFormat = "@"
$c.Value = "9"
$c.ClearFormats()
$d = $ws.Range("D97")
$d.NumberFormat = "@"
$d.Value = "21000.00"
$d.ClearFormats()

$c = $ws.Range("C98")
$c.NumberFormat = "@"
$c.Value = "38"
$c.ClearFormats()
$d = $ws.Range("D98")
$d.NumberFormat = "@"
$d.Value = "102500.00"
$d.ClearFormats()

$c = $ws.Range("C99")
$c.NumberFormat = "@"
$c.Value = "57"
$c.ClearFormats()
$d = $ws.Range("D99")
$d.NumberFormat = "@"
$d.Value = "139979.00"
$d.ClearFormats()

$c = $ws.Range("C100")
$c.NumberFormat = "@"
$c.Value = "21"
$c.ClearFormats()
$d = $ws.Range("D100")
$d.NumberFormat = "@"
$d.Value = "47700.00"
$d.ClearFormats()

$c = $ws.Range("C101")
$c.NumberFormat = "@"
$c.Value = "12"
$c.ClearFormats()
$d = $ws.Range("D101")
$d.NumberFormat = "@"
$d.Value = "30000.00"
$d.ClearFormats()

$c = $ws.Range("C102")
$c.NumberFormat = "@"
$c.Value = "18"
$c.ClearFormats()
$d = $ws.Range("D102")
$d.NumberFormat = "@"
$d.Value = "46830.00"
$d.ClearFormats()

$c = $ws.Range("C103")
$c.NumberFormat = "@"
$c.Value = "46"
$c.ClearFormats()
$d = $ws.Range("D103")
$d.NumberFormat = "@"
$d.Value = "93500.00"
$d.ClearFormats()

$c = $ws.Range("C104")
$c.NumberFormat = "@"
$c.Value = "7"
$c.ClearFormats()
$d = $ws.Range("D104")
$d.NumberFormat = "@"
$d.Value = "21909.00"
$d.ClearFormats()

$c = $ws.Range("C106")
$c.NumberFormat = "@"
$c.Value = "41"
$c.ClearFormats()
$d = $ws.Range("D106")
$d.NumberFormat = "@"
$d.Value = "107360.00"
$d.ClearFormats()

$c = $ws.Range("C108")
$c.NumberFormat = "@"
$c.Value = "8"
$c.ClearFormats()
$d = $ws.Range("D108")
$d.NumberFormat = "@"
$d.Value = "41219.00"
$d.ClearFormats()

$c = $ws.Range("C109")
$c.NumberFormat = "@"
$c.Value = "58"
$c.ClearFormats()
$d = $ws.Range("D109")
$d.NumberFormat = "@"
$d.Value = "348657.15"
$d.ClearFormats()

$c = $ws.Range("C110")
$c.NumberFormat = "@"
$c.Value = "5"
$c.ClearFormats()
$d = $ws.Range("D110")
$d.NumberFormat = "@"
$d.Value = "11500.00"
$d.ClearFormats()

$c = $ws.Range("C113")
$c.NumberFormat = "@"
$c.Value = "19"
$c.ClearFormats()
$d = $ws.Range("D113")
$d.NumberFormat = "@"
$d.Value = "49895.00"
$d.ClearFormats()

$c = $ws.Range("C115")
$c.NumberFormat = "@"
$c.Value = "6"
$c.ClearFormats()
$d = $ws.Range("D115")
$d.NumberFormat = "@"
$d.Value = "13500.00"
$d.ClearFormats()

$c = $ws.Range("C117")
$c.NumberFormat = "@"
$c.Value = "18"
$c.ClearFormats()
$d = $ws.Range("D117")
$d.NumberFormat = "@"
$d.Value = "65666.00"
$d.ClearFormats()

$c = $ws.Range("C122")
$c.NumberFormat = "@"
$c.Value = "50"
$c.ClearFormats()
$d = $ws.Range("D122")
$d.NumberFormat = "@"
$d.Value = "144676.45"
$d.ClearFormats()
